$d = $word.ActiveDocument

# 1. Switch the licence text from "CC BY-NC 4.0" to "CC BY-SA 4.0"
$d.Content.Find.Execute(
    "licensed under CC BY-NC 4.0. To view",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "licensed under CC BY-SA 4.0. To view", 2
)

# 2. Re-stamp the italics on just the "CC BY-SA" run so it becomes its own
#    run (matching how the licence name is now called out separately from
#    the surrounding sentence), even though the visual formatting does not
#    change.
$ccRange = $d.Content
$found = $ccRange.Find.Execute("CC BY-SA", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ccRange.Italic = 0
    $ccRange.Italic = 1
}

# 3. Update the licence hyperlink: address + display text, by-nc -> by-sa
$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.Address = "https://creativecommons.org/licenses/by-sa/4.0"
$hyperlink.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
